$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear B2 and B3 (previously 'wl_ com' and 'wl_dt')
$ws.Range("B2").Value = ""
$ws.Range("B3").Value = ""

# Update C2 and C3 with the new long text content
$ws.Range("C2").Value = 'Progress 0/ 13 00 mins : 00 secs 0.83 Understanding 0.94 Expressing Leaderboard Info Sign out Copyright Leaderboard Your Communication Score To know more about your performance in the simulation, download your personalized report. DOWNLOAD REPORT'
$ws.Range("C3").Value = '00 mins : 00 secs Viz Star 0.09 Insight-o-Meter 1.37% Leaderboard Info Info  Conversations 2 Sent Filter by: 5 Conversations Ambarjeet Kaur Director, Sales Sales Trends - All Channels - Good morning, Ax-Ax Could you turn last year''s monthly sales figures for all our channels (Direct, Retail and Online) into a nice, clean graph? I want to compare these channels'' performances individually and against each other - maybe there are months when we should focus on one channel more than the others.. Thanks. Ambarjeet 3 mins ago Alina Graf Director, Products Pending Qk Qn! Which one''s best? - Ax-Ax Hey! I need to pick your brains again. There are slides where I show trends across a specific timeframe - like, say, sales in a particular region over the past few months. I need to be able to read, at a glance, if the numbers are going up or down. What would you recommend I do? Regs, A. Graf 35 mins ago Alina Graf Director, Products Direct message from Alina - View the complete transcript 35 mins ago Ambarjeet Kaur Director, Sales Welcome to the Task Force! - Good morning, Alexa Hope you''ve had a chance to settle in and familiarize yourself with our organization, our team and our goals. And, as you may have heard, we''ve been going through a tough time of late. Business is down, and the CEO, Pierre Chartreuse, has set up a task force to figure out what has gone wrong and what we need to fix. I hope you don''t mind that I co-opted you into this task force. I''ve been quite impressed with the infographics you sent out as part of your application, and it''s exactly the kind of skillset we are missing at the moment. Think of it this way - it''s a good way to get noticed by the leadership team as quickly as possible. Your dashboard has a couple of indicators at the top that you may have noticed - these are part of our performance appraisal system. At TruFashione, we make it a point to be fair and transparent in recognizing the value that each and every person brings to work. As an expert in data visualization, you will be rated on the VizStar and Insight-o-Meter scales based on feedback from other members of the Task Force. Once again, welcome to TruFashione! Ambarjeet PS: I heard that your nickname is Ax-Ax. That''s a cool nickname! 40 mins ago Minnie Maxwell Director, HR Welcome to TruFashione! - Hi Alexa!  I hope you''ve had a seamless onboarding experience so far. We are quite pleased to have you with us, and you couldn''t have joined at a more exciting time!  This Inbox is where you''ll see all the official messages that come your way. That means emails, chats, meetings... if it''s digital, it happens here. Since we operate in a decentralized manner, this inbox can help you refer back to past conversations. But do remember that there is no RECALL facility - words once exchanged, unlike our goods (ha-ha!), will not be taken back under any circumstances.   Once again, welcome. Hope we have a great time together!  Warm regards, MM, Director, HR 40 mins ago 2024©KNOLSKAPE | All rights reserved Leaderboard Data Visualization Score To know more about your performance in the simulation, download your personalized report. Download Report'

# Widen column C to fit the new, much longer content
$ws.Columns.Item(3).ColumnWidth = 235.6

# Update the active selection to C4
$ws.Range("C4").Select()
